# Generate Report for Handback
#
# The handback pipeline finished: both localized files (zh-cn and de-de) for
# the two source docs are now in sync with en-US. Update the per-language
# status sheets with the resolved target/handback file links + handback
# timestamps, flip the status text on every sheet, and widen the columns
# that now hold longer content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/119d3efff385da3e2fd6917109bd721644865f4c/e2e/4298d119-7fdc-47de-9418-d452044444e2.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/119d3efff385da3e2fd6917109bd721644865f4c/e2e/ebc5170b-159a-47ae-b605-5af30015c4aa.md"
$nameA = "4298d119-7fdc-47de-9418-d452044444e2.md"
$nameB = "ebc5170b-159a-47ae-b605-5af30015c4aa.md"

# Excel's COM ColumnWidth setter adds back the standard ~5px padding
# (5/6 of a character unit) once the value round-trips through xlsx; undo
# that offset up front so the stored <col width=.../> lands on target.
$padding = 5 / 6
$wideWidth = 29.9777047293527 - $padding
$maxWidth = 40 - $padding

# ---------------------------------------------------------------------
# Overview sheet: status text propagates to the zh-cn/de-de summary
# columns; those two columns also get wider to fit the new text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("E1").ColumnWidth = $wideWidth
$wsOverview.Range("F1").ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, "", "", $nameA)
$wsZh.Range("J2").Value = "4298d119-7fdc-47de-9418-d452044444e2.fb9a4a9ec1ea1cc3425b98af3fd3196e03a0d60e.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-07 06:42:34"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlB, "", "", $nameB)
$wsZh.Range("J3").Value = "ebc5170b-159a-47ae-b605-5af30015c4aa.d4a0cf8123f5d086c585116d6691a719f2191c69.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-07 06:42:34"

$wsZh.Range("C1").ColumnWidth = $wideWidth
$wsZh.Range("I1").ColumnWidth = $maxWidth
$wsZh.Range("J1").ColumnWidth = $maxWidth

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, "", "", $nameA)
$wsDe.Range("J2").Value = "4298d119-7fdc-47de-9418-d452044444e2.fb9a4a9ec1ea1cc3425b98af3fd3196e03a0d60e.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-07 06:42:43"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlB, "", "", $nameB)
$wsDe.Range("J3").Value = "ebc5170b-159a-47ae-b605-5af30015c4aa.d4a0cf8123f5d086c585116d6691a719f2191c69.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-07 06:42:43"

$wsDe.Range("C1").ColumnWidth = $wideWidth
$wsDe.Range("I1").ColumnWidth = $maxWidth
$wsDe.Range("J1").ColumnWidth = $maxWidth
